# This script updates the "想去人数" (want-to-go count) column F values
# for several conventions, reflecting refreshed scrape data (gh-pages
# regeneration), on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1823
$wsExpo.Range("F6").Value = 675
$wsExpo.Range("F12").Value = 86
$wsExpo.Range("F13").Value = 173
$wsExpo.Range("F18").Value = 5212
$wsExpo.Range("F20").Value = 853
$wsExpo.Range("F21").Value = 124
$wsExpo.Range("F22").Value = 2307
$wsExpo.Range("F25").Value = 2152

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1823
$wsAll.Range("F6").Value = 675
$wsAll.Range("F12").Value = 86
$wsAll.Range("F13").Value = 173
$wsAll.Range("F18").Value = 5212
$wsAll.Range("F22").Value = 853
$wsAll.Range("F23").Value = 124
$wsAll.Range("F24").Value = 2307
$wsAll.Range("F28").Value = 2152
